$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns stay as text, matching the source data which
# is stored as inline strings (e.g. "27.590.42" is not a valid number).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.590.42"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "1.862.39"
$ws.Range("E3").Value = "  -4.46%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -1.05%  "
$ws.Range("D5").Value = "324.66"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "0.4514"
$ws.Range("E7").Value = "  -5.05%  "
$ws.Range("D8").Value = "0.3877"
$ws.Range("E8").Value = "  -3.75%  "
$ws.Range("D9").Value = "48.51"
$ws.Range("E9").Value = "  -9.61%  "
$ws.Range("D10").Value = "0.08043"
$ws.Range("E10").Value = "  -5.11%  "
$ws.Range("D11").Value = "1.022"
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("D12").Value = "21.63"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "1.907.62"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "7.179"
$ws.Range("E14").Value = "  -5.52%  "
$ws.Range("D15").Value = "5.888"
$ws.Range("E15").Value = "  -4.97%  "
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "0.00001040"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "86.25"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").Value = "0.06553"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "17.20"
$ws.Range("E20").Value = "  -7.56%  "
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "5.528"
$ws.Range("E22").Value = "  -4.76%  "
$ws.Range("D23").Value = "27.632.03"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  -5.30%  "
$ws.Range("D25").Value = "2.295"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "2.131.78"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").Value = "151.51"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "19.52"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").Value = "5.542"
$ws.Range("E29").Value = "  -6.37%  "
$ws.Range("D30").Value = "2.033"
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("D31").Value = "120.89"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").Value = "0.09404"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").Value = "1.469"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").Value = "0.9305"
$ws.Range("E34").Value = "  -6.33%  "
$ws.Range("D35").Value = "3.642"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "5.305"
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("D39").Value = "0.05998"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").Value = "8.422"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").Value = "0.5975"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").Value = "0.1863"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").Value = "10.33"
$ws.Range("E44").Value = "  -6.61%  "
$ws.Range("D45").Value = "1.281"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("D46").Value = "0.5682"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").Value = "12.29"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("D50").Value = "0.06864"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "1.007"
$ws.Range("E51").Value = "  -0.86%  "

# Row 37/38 swap (VeChain <-> TrustWalletToken)
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.231"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02237"
$ws.Range("E38").Value = "  -3.96%  "

# Row 48/49 swap (PancakeSwap <-> NEARProtocol)
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.949"
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "3.387"
$ws.Range("E49").Value = "  -0.54%  "
